# Add a new "Supplier" import column to the complex product import template.
# This mirrors the commit "Add data import permissions": a new column K is
# added so that a supplier can be specified for each imported product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column K.
$ws.Range("K1").Value = "Supplier"

# Give every existing data row (2-7) a supplier id of 2.
$ws.Range("K2").Value = 2
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 2

# Leave the selection on the newly-added cell, as in the source workbook.
$ws.Range("K7").Select()
